$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5225439071655273
$ws.Range("B1").Value = 1.522161960601807
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.470676183700562
$ws.Range("E1").Value = 1.417742848396301
